$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values (e.g. "1.003")
# are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '28.010.27', '  -0.30%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.863.91', '  -0.82%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.003', '  -0.40%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '312.19', '  -0.33%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  -0.27%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5128', '  +2.30%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3881', '  +1.56%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08354', '  -2.33%  '),
    @(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.110', '  -0.43%  '),
    @(11, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '41.37', '  -0.36%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.917.25', '  +1.78%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.156', '  -1.99%  '),
    @(14, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.42', '  -0.68%  '),
    @(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.256', '  +0.66%  '),
    @(16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  -0.50%  '),
    @(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001096', '  -0.30%  '),
    @(18, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '90.57', '  -0.76%  '),
    @(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06630', '  -0.06%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.61', '  -2.50%  '),
    @(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  -0.16%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.993', '  -1.55%  '),
    @(23, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '28.054.62', '  -0.32%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.02', '  -2.05%  '),
    @(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.241', '  -1.88%  '),
    @(26, 'LEO', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', '3.399', '  -0.59%  '),
    @(27, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.074.09', '  -1.22%  '),
    @(28, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '158.09', '  +0.94%  '),
    @(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.459', '  -4.90%  '),
    @(30, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '20.49', '  -1.05%  '),
    @(31, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '124.38', '  -1.36%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1061', '  +0.59%  '),
    @(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.027', '  -2.22%  '),
    @(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.836', '  +3.64%  '),
    @(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.600', '  -0.15%  '),
    @(36, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.378', '  -2.11%  '),
    @(37, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02431', '  -0.95%  '),
    @(38, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06528', '  +0.04%  '),
    @(39, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2176', '  +0.27%  '),
    @(40, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.196', '  -2.98%  '),
    @(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6480', '  +1.89%  '),
    @(42, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.970', '  +1.86%  '),
    @(43, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.29', '  -0.58%  '),
    @(44, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.209', '  -2.18%  '),
    @(45, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6058', '  +0.42%  '),
    @(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.97', '  -1.07%  '),
    @(47, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.288', '  -0.94%  '),
    @(48, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.672', '  -0.21%  '),
    @(49, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.995', '  -0.15%  '),
    @(50, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.218', '  -0.11%  '),
    @(51, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '120.70', '  +0.01%  ')
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
    $ws.Cells.Item($r, 5).Value = $entry[4]
}
